## Applies the "EoCEDwEC" elasticity-source update:
##   - "About" sheet: refresh source citation (new EIA paper, 2014, Appendix)
##     and update the "1-year" -> "3-year" short-run-elasticity note text.
##   - "EIA Table 1" sheet: replace the short/long-run elasticity figures
##     with the values from the new source.
##   - "EoCEDwEC" sheet: the long-run-minus-short-run formulas now subtract
##     the 3-year (column D) short-run figure instead of the 1-year
##     (column B) figure.
##   - Restore sheet selections / active tab to match the saved workbook
##     state (EoCEDwEC becomes the active tab).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "About" sheet
# ---------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("B4").Value = 2014
$wsAbout.Range("B5").Value = "Price Elasticities for Energy Use in Buildings of the United States"
$wsAbout.Range("B6").Value = "https://www.eia.gov/analysis/studies/buildings/energyuse/pdf/price_elasticities.pdf"
$wsAbout.Range("B7").Value = "Appendix"

$wsAbout.Range("A10").Value = "We use same-price, long-run elasticities minus the 3-year short-run elasticities."
$wsAbout.Range("A11").Value = "We calculate it this way because we assume that 3-year elasticities primarily reflect behavior"
$wsAbout.Range("A14").Value = "all timescales.  So, the portion of the long-run elasticitiy represented by the 3-year elasticity"

# ---------------------------------------------------------------------
# "EIA Table 1" sheet
# ---------------------------------------------------------------------
$wsEia = $wb.Worksheets.Item("EIA Table 1")

# Residential section (row 7 = Electricity, 8 = Natural Gas, 9 = Distillate Fuel)
$wsEia.Range("B7").Value = -0.12
$wsEia.Range("C7").Value = -0.21
$wsEia.Range("D7").Value = -0.25
$wsEia.Range("E7").Value = -0.28
$wsEia.Range("F7").Value = 0
$wsEia.Range("G7").Value = 0

$wsEia.Range("B8").Value = -0.07
$wsEia.Range("C8").Value = -0.13
$wsEia.Range("D8").Value = -0.15
$wsEia.Range("E8").Value = 0.03
$wsEia.Range("F8").Value = -0.21
$wsEia.Range("G8").Value = 0

$wsEia.Range("B9").Value = -0.07
$wsEia.Range("C9").Value = -0.12
$wsEia.Range("D9").Value = -0.14
$wsEia.Range("E9").Value = 0
$wsEia.Range("F9").Value = 0
$wsEia.Range("G9").Value = -0.22

# Commercial section (row 14 = Electricity, 15 = Natural Gas, 16 = Distillate Fuel)
$wsEia.Range("B14").Value = -0.11
$wsEia.Range("C14").Value = -0.18
$wsEia.Range("D14").Value = -0.22
$wsEia.Range("E14").Value = -0.33
$wsEia.Range("F14").Value = 0.09
$wsEia.Range("G14").Value = 0

$wsEia.Range("B15").Value = -0.15
$wsEia.Range("C15").Value = -0.25
$wsEia.Range("D15").Value = -0.3
$wsEia.Range("E15").Value = 0.15
$wsEia.Range("F15").Value = -0.58
$wsEia.Range("G15").Value = 0.02

$wsEia.Range("B16").Value = -0.14
$wsEia.Range("C16").Value = -0.24
$wsEia.Range("D16").Value = -0.29
$wsEia.Range("E16").Value = 0
$wsEia.Range("F16").Value = 0.05
$wsEia.Range("G16").Value = -0.42

# ---------------------------------------------------------------------
# "EoCEDwEC" sheet — switch the long-run-minus-short-run formulas from
# the 1-year column (B) to the 3-year column (D) on "EIA Table 1"
# ---------------------------------------------------------------------
$wsMain = $wb.Worksheets.Item("EoCEDwEC")

$wsMain.Range("B2").Formula = "='EIA Table 1'!E7-'EIA Table 1'!D7"
$wsMain.Range("D2").Formula = "='EIA Table 1'!E14-'EIA Table 1'!D14"
$wsMain.Range("B4").Formula = "='EIA Table 1'!F8-'EIA Table 1'!D8"
$wsMain.Range("D4").Formula = "='EIA Table 1'!F15-'EIA Table 1'!D15"
$wsMain.Range("B5").Formula = "='EIA Table 1'!G9-'EIA Table 1'!D9"
$wsMain.Range("D5").Formula = "='EIA Table 1'!G16-'EIA Table 1'!D16"

# ---------------------------------------------------------------------
# Selections / active tab — EoCEDwEC ends up the active sheet, with the
# other two sheets remembering a last-used selection.
# ---------------------------------------------------------------------
$wsAbout.Activate()
$wsAbout.Range("A27").Select()

$wsEia.Activate()
$wsEia.Range("E17").Select()

$wsMain.Activate()
$wsMain.Range("H29").Select()
